$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 190, shifting existing rows 190:246 down to 191:247
$ws.Rows.Item(190).Insert()

# Populate the newly inserted row 190 with the new data record
$ws.Range("A190").Value2 = 10
$ws.Range("B190").Value2 = "Vega Modelo de Temuco"
$ws.Range("C190").Value2 = "La Araucanía"
$ws.Range("D190").Value2 = 44627
$ws.Range("E190").Value2 = 9
$ws.Range("F190").Value2 = 100112001
$ws.Range("G190").Value2 = "Berenjena"
$ws.Range("H190").Value2 = "Sin especificar"
$ws.Range("I190").Value2 = "Primera"
$ws.Range("J190").Value2 = 50
$ws.Range("K190").Value2 = 10000
$ws.Range("L190").Value2 = 10000
$ws.Range("M190").Value2 = 10000
$ws.Range("N190").Value2 = "$/caja 60 unidades"
$ws.Range("O190").Value2 = "Región del Maule"
$ws.Range("P190").Value2 = 167
$ws.Range("Q190").Value2 = 60
$ws.Range("R190").Value2 = "Hortaliza"
